$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2020
$ws.Range("D2").Value = 11278

$ws.Range("B4").Value = 10100
$ws.Range("C4").Value = 9400
$ws.Range("D4").Value = 29940
$ws.Range("E4").Value = 12480
$ws.Range("F4").Value = 28703

$ws.Range("D9").Select()
